$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Geeks-for-Geeks row (row 2): C2 13->14, D2 24->25 (F2 recalculates via formula)
$ws.Range("C2").Value = 14
$ws.Range("D2").Value = 25

# Update Leetcode row (row 4): E4 16->20 (F4 recalculates via formula)
$ws.Range("E4").Value = 20

# Update the active selection to E5
$ws.Range("E5").Select()

$wb.Save()
